$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite the existing "url" column (A2:A30) from the old
# "YYYYMMDDNNNNSUFFIX" format into "YYYY-MM-DD_REG-NNN_SUFFIX", where NNN
# is the row's "idingreso" (column B) value, zero padded to 3 digits.
for ($r = 2; $r -le 30; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $url = $cell.Value2
    $idingreso = $ws.Cells.Item($r, 2).Value2

    if ($url -match '^(\d{4})(\d{2})(\d{2})\d{4}([A-Za-z]+)$') {
        $year = $matches[1]
        $month = $matches[2]
        $day = $matches[3]
        $suffix = $matches[4]
        $reg = "{0:D3}" -f [int]$idingreso
        $cell.Value = "$year-$month-$day`_REG-$reg`_$suffix"
    }
}

# Append the new row (31) for the newly loaded "ingreso" record.
$ws.Cells.Item(31, 1).Value = "2017-11-07_REG-025_BOLEC"
$ws.Cells.Item(31, 2).Value = 25
$ws.Cells.Item(31, 3).Value = 22

# Move the selection to the new row, matching where the user ended up
# after typing the new record.
$ws.Range("A31").Select()
